$d = $word.ActiveDocument

# The document has a single paragraph whose only run holds just a space
# character:
#     <w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
#
# The target splits that into two runs with identical run properties:
#     <w:r>...<w:t>Rio is best</w:t></w:r>
#     <w:r>...<w:t xml:space="preserve"> in everything </w:t></w:r>
#
# A plain Range.Text assignment (or InsertAfter) would leave two adjacent
# runs with identical formatting, which get coalesced into a single run on
# save. Using Range.InsertXML with an explicit WordprocessingML fragment
# preserves the run boundary exactly as authored.

$para = $d.Paragraphs.Item(1)
$full = $para.Range

# Range covering the paragraph's text but not its paragraph mark, so the
# inserted XML replaces just the run content and keeps the existing <w:p>
# (with its pPr/paraId/etc.) untouched.
$target = $d.Range($full.Start, $full.End - 1)

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-IN"/>
              </w:rPr>
              <w:t>Rio is best</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-IN"/>
              </w:rPr>
              <w:t xml:space="preserve"> in everything </w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$target.InsertXML($xml)

$d.Save()
